$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds plain-looking decimal strings (e.g. "7.20")
# that Excel would normally auto-coerce to a Number, which silently drops
# trailing zeros (e.g. "7.20" -> 7.2). Force text storage, then clear the
# temporary number-format override so the cell keeps its original (default)
# style, matching the unaffected sibling cells.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextValue $ws.Range('D2') '64.965.20'
$ws.Range('E2').Value = '  +0.95%  '
Set-TextValue $ws.Range('D3') '3.186.88'
$ws.Range('E3').Value = '  +1.24%  '
$ws.Range('E4').Value = '  +0.18%  '
Set-TextValue $ws.Range('D5') '615.18'
$ws.Range('E5').Value = '  +1.53%  '
Set-TextValue $ws.Range('D6') '147.77'
$ws.Range('E6').Value = '  -1.17%  '
$ws.Range('E7').Value = '  +0.00%  '
Set-TextValue $ws.Range('D8') '3.181.27'
$ws.Range('E8').Value = '  +1.09%  '
$ws.Range('E10').Value = '  +0.37%  '
Set-TextValue $ws.Range('D11') '5.52'
$ws.Range('E11').Value = '  -1.29%  '
Set-TextValue $ws.Range('D12') '0.479'
$ws.Range('E12').Value = '  -0.07%  '
Set-TextValue $ws.Range('D13') '0.0000264'
$ws.Range('E13').Value = '  +1.54%  '
Set-TextValue $ws.Range('D14') '36.18'
$ws.Range('E14').Value = '  -1.85%  '
Set-TextValue $ws.Range('D15') '3.710.96'
$ws.Range('E15').Value = '  +1.41%  '
$ws.Range('E16').Value = '  +3.26%  '
Set-TextValue $ws.Range('D17') '65.002.62'
$ws.Range('E17').Value = '  +0.99%  '
Set-TextValue $ws.Range('D18') '3.189.94'
$ws.Range('E18').Value = '  +1.47%  '
Set-TextValue $ws.Range('D19') '6.96'
$ws.Range('E19').Value = '  -0.68%  '
Set-TextValue $ws.Range('D20') '484.82'
$ws.Range('E20').Value = '  +0.27%  '
Set-TextValue $ws.Range('D21') '14.81'
$ws.Range('E21').Value = '  +1.02%  '
$ws.Range('E22').Value = '  +1.98%  '
Set-TextValue $ws.Range('D23') '7.98'
$ws.Range('E23').Value = '  +2.73%  '
Set-TextValue $ws.Range('D24') '13.95'
$ws.Range('E24').Value = '  +0.77%  '
Set-TextValue $ws.Range('D25') '84.67'
$ws.Range('E25').Value = '  +0.64%  '
$ws.Range('E26').Value = '  -0.03%  '
Set-TextValue $ws.Range('D27') '8.85'
$ws.Range('E27').Value = '  +3.40%  '
$ws.Range('E28').Value = '  -3.51%  '
Set-TextValue $ws.Range('D29') '7.20'
$ws.Range('E29').Value = '  +3.48%  '
$ws.Range('E30').Value = '  -3.57%  '
$ws.Range('E31').Value = '  -4.92%  '
$ws.Range('E32').Value = '  +0.68%  '
$ws.Range('E33').Value = '  +0.21%  '
Set-TextValue $ws.Range('D34') '26.83'
$ws.Range('E34').Value = '  +0.20%  '
$ws.Range('E35').Value = '  +2.62%  '
Set-TextValue $ws.Range('D36') '0.0₃0799'
$ws.Range('E36').Value = '  +5.66%  '
Set-TextValue $ws.Range('D37') '6.06'
$ws.Range('E37').Value = '  -0.81%  '
Set-TextValue $ws.Range('D38') '3.22'
$ws.Range('E38').Value = '  -1.47%  '
Set-TextValue $ws.Range('D39') '53.39'
$ws.Range('E39').Value = '  -1.83%  '
Set-TextValue $ws.Range('D40') '470.15'
$ws.Range('E40').Value = '  +3.85%  '
$ws.Range('E41').Value = '  +0.59%  '
$ws.Range('E42').Value = '  -1.87%  '
Set-TextValue $ws.Range('D43') '8.43'
$ws.Range('E43').Value = '  -0.78%  '
Set-TextValue $ws.Range('D44') '2.873.77'
$ws.Range('E44').Value = '  -0.46%  '
$ws.Range('E45').Value = '  +1.95%  '
$ws.Range('E46').Value = '  -0.14%  '
Set-TextValue $ws.Range('D47') '2.47'
$ws.Range('E47').Value = '  +6.14%  '
Set-TextValue $ws.Range('D48') '37.32'
$ws.Range('E48').Value = '  +10.91%  '
Set-TextValue $ws.Range('D49') '27.03'
$ws.Range('E49').Value = '  +0.98%  '
$ws.Range('E50').Value = '  +0.09%  '
$ws.Range('E51').Value = '  -0.49%  '
